$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$ws.Range('D2').Value = '34.086.14'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = '1.786.87'
$ws.Range('E3').Value = '  -2.89%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'D5' '224.26'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('E6').Value = '  -1.43%  '
$ws.Range('E7').Value = '  +0.08%  '
Set-TextValue 'D8' '32.76'
$ws.Range('E8').Value = '  +2.30%  '
$ws.Range('E9').Value = '  -2.39%  '
$ws.Range('E10').Value = '  -1.11%  '
Set-TextValue 'D11' '0.0936'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '2.044.82'
$ws.Range('E12').Value = '  -2.83%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D13' '10.94'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.760.16'
$ws.Range('E14').Value = '  -4.83%  '
$ws.Range('E15').Value = '  -3.99%  '
$ws.Range('D16').Value = '34.044.15'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('E17').Value = '  -4.35%  '
$ws.Range('E18').Value = '  -3.03%  '
Set-TextValue 'D19' '245.50'
$ws.Range('E19').Value = '  -2.71%  '
$ws.Range('D20').Value = '0.0₃0790'
$ws.Range('E20').Value = '  -1.02%  '
Set-TextValue 'D22' '10.84'
$ws.Range('E22').Value = '  -4.64%  '
$ws.Range('E23').Value = '  -4.64%  '
$ws.Range('E24').Value = '  -3.13%  '
Set-TextValue 'D25' '160.67'
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('E27').Value = '  -2.89%  '
$ws.Range('E28').Value = '  -2.56%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('E32').Value = '  -4.16%  '
$ws.Range('E34').Value = '  -5.66%  '
$ws.Range('D35').Value = '1.396.72'
$ws.Range('E35').Value = '  -4.35%  '
Set-TextValue 'D36' '0.645'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('E37').Value = '  -1.40%  '
$ws.Range('E38').Value = '  -3.41%  '
$ws.Range('E39').Value = '  +2.52%  '
$ws.Range('E40').Value = '  -0.15%  '
Set-TextValue 'D41' '0.917'
$ws.Range('E41').Value = '  -5.49%  '
$ws.Range('E42').Value = '  -2.87%  '
Set-TextValue 'D43' '77.96'
$ws.Range('E43').Value = '  -5.31%  '
$ws.Range('D44').Value = '0.0₆0146'
$ws.Range('E44').Value = '  +16.14%  '
$ws.Range('E45').Value = '  +1.86%  '
Set-TextValue 'D46' '12.52'
$ws.Range('E46').Value = '  +4.12%  '
$ws.Range('E47').Value = '  -0.59%  '
Set-TextValue 'D48' '108.29'
$ws.Range('E48').Value = '  +1.61%  '
Set-TextValue 'D49' '5.87'
$ws.Range('E49').Value = '  -3.63%  '
$ws.Range('D50').Value = '1.944.54'
$ws.Range('E50').Value = '  -2.75%  '
$ws.Range('E51').Value = '  +0.05%  '
